$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster = M2, Ligand = Spn, Receptor = Siglec1, Target cluster = M2
$ws.Range("A2").Value = "M2"
$ws.Range("B2").Value = "Spn"
$ws.Range("C2").Value = "Siglec1"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.400925
$ws.Range("H2").Value = 7.202775
$ws.Range("I2").Value = 0.9592478079643895
$ws.Range("J2").Value = 0.9592478079643894
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 55.49088033333334
$ws.Range("N2").Value = 166.472641
$ws.Range("O2").Value = 0.9642568285787312
$ws.Range("P2").Value = 0.9642568285787311
$ws.Range("Q2").Value = 133.2294418643083
$ws.Range("R2").Value = 1199.064976778775
$ws.Range("S2").Value = 0.9249612491288419
$ws.Range("T2").Value = 0.9249612491288417

# Row 3: Sending cluster = M2, Ligand = Spn, Receptor = Siglec1, Target cluster = sCs
$ws.Range("A3").Value = "M2"
$ws.Range("B3").Value = "Spn"
$ws.Range("C3").Value = "Siglec1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.400925
$ws.Range("H3").Value = 7.202775
$ws.Range("I3").Value = 0.9592478079643895
$ws.Range("J3").Value = 0.9592478079643894
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 2.056941666666666
$ws.Range("N3").Value = 6.170825
$ws.Range("O3").Value = 0.03574317142126884
$ws.Range("P3").Value = 0.03574317142126884
$ws.Range("Q3").Value = 4.938562671041666
$ws.Range("R3").Value = 44.447064039375
$ws.Range("S3").Value = 0.03428655883554755
$ws.Range("T3").Value = 0.03428655883554754

# Row 4: Sending cluster = sCs, Ligand = Spn, Receptor = Siglec1, Target cluster = M2
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Spn"
$ws.Range("C4").Value = "Siglec1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1019996666666667
$ws.Range("H4").Value = 0.305999
$ws.Range("I4").Value = 0.04075219203561061
$ws.Range("J4").Value = 0.04075219203561061
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 55.49088033333334
$ws.Range("N4").Value = 166.472641
$ws.Range("O4").Value = 0.9642568285787312
$ws.Range("P4").Value = 0.9642568285787311
$ws.Range("Q4").Value = 5.66005129703989
$ws.Range("R4").Value = 50.940461673359
$ws.Range("S4").Value = 0.03929557944988932
$ws.Range("T4").Value = 0.03929557944988931

# Row 5: Sending cluster = sCs, Ligand = Spn, Receptor = Siglec1, Target cluster = sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Spn"
$ws.Range("C5").Value = "Siglec1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1019996666666667
$ws.Range("H5").Value = 0.305999
$ws.Range("I5").Value = 0.04075219203561061
$ws.Range("J5").Value = 0.04075219203561061
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 2.056941666666666
$ws.Range("N5").Value = 6.170825
$ws.Range("O5").Value = 0.03574317142126884
$ws.Range("P5").Value = 0.03574317142126884
$ws.Range("Q5").Value = 0.2098073643527777
$ws.Range("R5").Value = 1.888266279175
$ws.Range("S5").Value = 0.001456612585721297
$ws.Range("T5").Value = 0.001456612585721297
